$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$c = $ws.Range("AB1").AddComment("Hello World")
$c.Author = "Matt"
